$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 11457
$ws.Range("F3").Value = 1968
$ws.Range("F4").Value = 541
$ws.Range("F5").Value = 856
$ws.Range("F6").Value = 2435
$ws.Range("F7").Value = 776
$ws.Range("F9").Value = 612
$ws.Range("F10").Value = 471
$ws.Range("F11").Value = 1378
$ws.Range("F12").Value = 699
$ws.Range("F13").Value = 135
$ws.Range("F15").Value = 1008
$ws.Range("F16").Value = 562
$ws.Range("F17").Value = 691
$ws.Range("F18").Value = 1150
$ws.Range("F19").Value = 219
$ws.Range("F20").Value = 949
$ws.Range("F21").Value = 17
$ws.Range("F22").Value = 148
$ws.Range("F23").Value = 7
$ws.Range("F24").Value = 325
$ws.Range("F26").Value = 270
$ws.Range("F27").Value = 480
$ws.Range("F28").Value = 509
$ws.Range("F29").Value = 695
$ws.Range("F30").Value = 187
$ws.Range("F32").Value = 335

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 912
$ws.Range("F7").Value = 67
$ws.Range("F10").Value = 47
$ws.Range("F11").Value = 419

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 88

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 11457
$ws.Range("F3").Value = 1968
$ws.Range("F5").Value = 541
$ws.Range("F6").Value = 856
$ws.Range("F7").Value = 2435
$ws.Range("F8").Value = 776
$ws.Range("F11").Value = 612
$ws.Range("F12").Value = 471
$ws.Range("F13").Value = 88
$ws.Range("F14").Value = 1378
$ws.Range("F16").Value = 699
$ws.Range("F17").Value = 135
$ws.Range("F18").Value = 912
$ws.Range("F20").Value = 1008
$ws.Range("F21").Value = 562
$ws.Range("F22").Value = 691
$ws.Range("F23").Value = 1150
$ws.Range("F24").Value = 219
$ws.Range("F25").Value = 949
$ws.Range("F26").Value = 17
$ws.Range("F27").Value = 148
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 325
$ws.Range("F31").Value = 67
$ws.Range("F33").Value = 270
$ws.Range("F37").Value = 480
$ws.Range("F38").Value = 509
$ws.Range("F39").Value = 695
$ws.Range("F40").Value = 187
$ws.Range("F41").Value = 47
$ws.Range("F43").Value = 419
$ws.Range("F45").Value = 335
